# Insert a new weekly price-report row at row 92 (shifting existing rows 92-175 down to 93-176)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(92).Insert()

$ws.Cells.Item(92, 1).Value  = 3
$ws.Cells.Item(92, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(92, 3).Value  = "Coquimbo"
$ws.Cells.Item(92, 4).Value  = 44658
$ws.Cells.Item(92, 5).Value  = 5
$ws.Cells.Item(92, 6).Value  = 100112030
$ws.Cells.Item(92, 7).Value  = "Poroto granado"
$ws.Cells.Item(92, 8).Value  = "Sin especificar"
$ws.Cells.Item(92, 9).Value  = "Primera"
$ws.Cells.Item(92, 10).Value = 73
$ws.Cells.Item(92, 11).Value = 21000
$ws.Cells.Item(92, 12).Value = 22000
$ws.Cells.Item(92, 13).Value = 21521
$ws.Cells.Item(92, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(92, 15).Value = "Provincia de Talca"
$ws.Cells.Item(92, 16).Value = 861
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
